$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new test case row (row 29), mirroring the formatting of the last
# existing row (row 28: test_CreateUser) so the "Outcome"/"Justification"
# columns (E:F) keep their usual applied style.
$ws.Rows("28:28").Copy()
$ws.Rows("29:29").Insert(-4121)

$ws.Range("A29").Value = 7
$ws.Range("B29").Value = "test_DeleteUser"
$ws.Range("C29").Value = "This is to test whether users are able to delete an existing user"

# D28 used a wrap-text style (long justification text); this new row's
# "Test Values" is just "NIL", so drop the inherited wrap formatting.
$ws.Range("D29").ClearFormats()
$ws.Range("D29").Value = "NIL"

$ws.Range("E29").Value = "User is deleteds"
$ws.Range("F29").Value = "Case failed"

# No "Justification" note for this new case.
$ws.Range("G29").ClearContents()

$ws.Range("G29").Select()
